$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-11 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-12 Friday", 2) | Out-Null
$d.Content.Find.Execute("63-45=18", $true, $false, $false, $false, $false, $true, 1, $false, "13+69=82", 2) | Out-Null
$d.Content.Find.Execute("30-15=15", $true, $false, $false, $false, $false, $true, 1, $false, "25-12=13", 2) | Out-Null
$d.Content.Find.Execute("97-65=32", $true, $false, $false, $false, $false, $true, 1, $false, "66-4=62", 2) | Out-Null
$d.Content.Find.Execute("94-29=65", $true, $false, $false, $false, $false, $true, 1, $false, "94-48=46", 2) | Out-Null
$d.Content.Find.Execute("21+42=63", $true, $false, $false, $false, $false, $true, 1, $false, "0+93=93", 2) | Out-Null
$d.Content.Find.Execute("86-70=16", $true, $false, $false, $false, $false, $true, 1, $false, "9+5=14", 2) | Out-Null
$d.Content.Find.Execute("57-26=31", $true, $false, $false, $false, $false, $true, 1, $false, "76-41=35", 2) | Out-Null
$d.Content.Find.Execute("10+29=39", $true, $false, $false, $false, $false, $true, 1, $false, "11+83=94", 2) | Out-Null
$d.Content.Find.Execute("87-14=73", $true, $false, $false, $false, $false, $true, 1, $false, "66-26=40", 2) | Out-Null
$d.Content.Find.Execute("82-14=68", $true, $false, $false, $false, $false, $true, 1, $false, "93-31=62", 2) | Out-Null
$d.Content.Find.Execute("70-9=61", $true, $false, $false, $false, $false, $true, 1, $false, "30+62=92", 2) | Out-Null
$d.Content.Find.Execute("97-61=36", $true, $false, $false, $false, $false, $true, 1, $false, "65-56=9", 2) | Out-Null
$d.Content.Find.Execute("0+32=32", $true, $false, $false, $false, $false, $true, 1, $false, "88-62=26", 2) | Out-Null
$d.Content.Find.Execute("48+18=66", $true, $false, $false, $false, $false, $true, 1, $false, "52+33=85", 2) | Out-Null
$d.Content.Find.Execute("54-30=24", $true, $false, $false, $false, $false, $true, 1, $false, "12+12=24", 2) | Out-Null
$d.Content.Find.Execute("39+55=94", $true, $false, $false, $false, $false, $true, 1, $false, "48-4=44", 2) | Out-Null
$d.Content.Find.Execute("70-39=31", $true, $false, $false, $false, $false, $true, 1, $false, "59+13=72", 2) | Out-Null
$d.Content.Find.Execute("17-7=10", $true, $false, $false, $false, $false, $true, 1, $false, "66-47=19", 2) | Out-Null
$d.Content.Find.Execute("10+59=69", $true, $false, $false, $false, $false, $true, 1, $false, "87+5=92", 2) | Out-Null
$d.Content.Find.Execute("6+14=20", $true, $false, $false, $false, $false, $true, 1, $false, "51+37=88", 2) | Out-Null
$d.Content.Find.Execute("30+55=85", $true, $false, $false, $false, $false, $true, 1, $false, "15+18=33", 2) | Out-Null
$d.Content.Find.Execute("10+1=11", $true, $false, $false, $false, $false, $true, 1, $false, "69-7=62", 2) | Out-Null
$d.Content.Find.Execute("44-40=4", $true, $false, $false, $false, $false, $true, 1, $false, "61+38=99", 2) | Out-Null
$d.Content.Find.Execute("36+44=80", $true, $false, $false, $false, $false, $true, 1, $false, "86-59=27", 2) | Out-Null
$d.Content.Find.Execute("1+25=26", $true, $false, $false, $false, $false, $true, 1, $false, "28+34=62", 2) | Out-Null
$d.Content.Find.Execute("30-11=19", $true, $false, $false, $false, $false, $true, 1, $false, "9+46=55", 2) | Out-Null
$d.Content.Find.Execute("36+7=43", $true, $false, $false, $false, $false, $true, 1, $false, "75-71=4", 2) | Out-Null
$d.Content.Find.Execute("99-39=60", $true, $false, $false, $false, $false, $true, 1, $false, "39-25=14", 2) | Out-Null
$d.Content.Find.Execute("98-58=40", $true, $false, $false, $false, $false, $true, 1, $false, "26+12=38", 2) | Out-Null
$d.Content.Find.Execute("80-4=76", $true, $false, $false, $false, $false, $true, 1, $false, "13-1=12", 2) | Out-Null
$d.Content.Find.Execute("85-41=44", $true, $false, $false, $false, $false, $true, 1, $false, "16+30=46", 2) | Out-Null
$d.Content.Find.Execute("17+58=75", $true, $false, $false, $false, $false, $true, 1, $false, "32-9=23", 2) | Out-Null
$d.Content.Find.Execute("42+55=97", $true, $false, $false, $false, $false, $true, 1, $false, "49-13=36", 2) | Out-Null
$d.Content.Find.Execute("91-45=46", $true, $false, $false, $false, $false, $true, 1, $false, "14+13=27", 2) | Out-Null
$d.Content.Find.Execute("10+69=79", $true, $false, $false, $false, $false, $true, 1, $false, "96-88=8", 2) | Out-Null
$d.Content.Find.Execute("2+96=98", $true, $false, $false, $false, $false, $true, 1, $false, "51-9=42", 2) | Out-Null
$d.Content.Find.Execute("94-47=47", $true, $false, $false, $false, $false, $true, 1, $false, "22+64=86", 2) | Out-Null
$d.Content.Find.Execute("66-65=1", $true, $false, $false, $false, $false, $true, 1, $false, "30+59=89", 2) | Out-Null
$d.Content.Find.Execute("3+92=95", $true, $false, $false, $false, $false, $true, 1, $false, "51-48=3", 2) | Out-Null
$d.Content.Find.Execute("9+16=25", $true, $false, $false, $false, $false, $true, 1, $false, "86-32=54", 2) | Out-Null
$d.Content.Find.Execute("54+4=58", $true, $false, $false, $false, $false, $true, 1, $false, "73-21=52", 2) | Out-Null
$d.Content.Find.Execute("61+4=65", $true, $false, $false, $false, $false, $true, 1, $false, "47-37=10", 2) | Out-Null
$d.Content.Find.Execute("97-27=70", $true, $false, $false, $false, $false, $true, 1, $false, "26+19=45", 2) | Out-Null
$d.Content.Find.Execute("47-19=28", $true, $false, $false, $false, $false, $true, 1, $false, "72+22=94", 2) | Out-Null
$d.Content.Find.Execute("30+1=31", $true, $false, $false, $false, $false, $true, 1, $false, "36+28=64", 2) | Out-Null
$d.Content.Find.Execute("50+15=65", $true, $false, $false, $false, $false, $true, 1, $false, "63-3=60", 2) | Out-Null
$d.Content.Find.Execute("5+27=32", $true, $false, $false, $false, $false, $true, 1, $false, "1+46=47", 2) | Out-Null
$d.Content.Find.Execute("26+72=98", $true, $false, $false, $false, $false, $true, 1, $false, "42+17=59", 2) | Out-Null
$d.Content.Find.Execute("6+2=8", $true, $false, $false, $false, $false, $true, 1, $false, "28-18=10", 2) | Out-Null
$d.Content.Find.Execute("4+84=88", $true, $false, $false, $false, $false, $true, 1, $false, "59+22=81", 2) | Out-Null
$d.Content.Find.Execute("7+44=51", $true, $false, $false, $false, $false, $true, 1, $false, "32-16=16", 2) | Out-Null
$d.Content.Find.Execute("22+77=99", $true, $false, $false, $false, $false, $true, 1, $false, "8+74=82", 2) | Out-Null
$d.Content.Find.Execute("96-45=51", $true, $false, $false, $false, $false, $true, 1, $false, "73-67=6", 2) | Out-Null
$d.Content.Find.Execute("96-18=78", $true, $false, $false, $false, $false, $true, 1, $false, "47-3=44", 2) | Out-Null
$d.Content.Find.Execute("76-15=61", $true, $false, $false, $false, $false, $true, 1, $false, "24+1=25", 2) | Out-Null
$d.Content.Find.Execute("1+48=49", $true, $false, $false, $false, $false, $true, 1, $false, "76-5=71", 2) | Out-Null
$d.Content.Find.Execute("23-18=5", $true, $false, $false, $false, $false, $true, 1, $false, "85-63=22", 2) | Out-Null
$d.Content.Find.Execute("83-23=60", $true, $false, $false, $false, $false, $true, 1, $false, "24+18=42", 2) | Out-Null
$d.Content.Find.Execute("66-10=56", $true, $false, $false, $false, $false, $true, 1, $false, "8+49=57", 2) | Out-Null
$d.Content.Find.Execute("10-6=4", $true, $false, $false, $false, $false, $true, 1, $false, "60+36=96", 2) | Out-Null
$d.Content.Find.Execute("58+18=76", $true, $false, $false, $false, $false, $true, 1, $false, "59+23=82", 2) | Out-Null
$d.Content.Find.Execute("85-1=84", $true, $false, $false, $false, $false, $true, 1, $false, "67+3=70", 2) | Out-Null
$d.Content.Find.Execute("8+73=81", $true, $false, $false, $false, $false, $true, 1, $false, "32+42=74", 2) | Out-Null
$d.Content.Find.Execute("65+31=96", $true, $false, $false, $false, $false, $true, 1, $false, "99-7=92", 2) | Out-Null
$d.Content.Find.Execute("2+15=17", $true, $false, $false, $false, $false, $true, 1, $false, "53-25=28", 2) | Out-Null
$d.Content.Find.Execute("28+66=94", $true, $false, $false, $false, $false, $true, 1, $false, "14+3=17", 2) | Out-Null
$d.Content.Find.Execute("78-23=55", $true, $false, $false, $false, $false, $true, 1, $false, "18+75=93", 2) | Out-Null
$d.Content.Find.Execute("12+22=34", $true, $false, $false, $false, $false, $true, 1, $false, "62-59=3", 2) | Out-Null
$d.Content.Find.Execute("92-6=86", $true, $false, $false, $false, $false, $true, 1, $false, "12+70=82", 2) | Out-Null
$d.Content.Find.Execute("68-34=34", $true, $false, $false, $false, $false, $true, 1, $false, "20+13=33", 2) | Out-Null
$d.Content.Find.Execute("37+27=64", $true, $false, $false, $false, $false, $true, 1, $false, "22+8=30", 2) | Out-Null
$d.Content.Find.Execute("86-68=18", $true, $false, $false, $false, $false, $true, 1, $false, "17+34=51", 2) | Out-Null
$d.Content.Find.Execute("36+52=88", $true, $false, $false, $false, $false, $true, 1, $false, "12+13=25", 2) | Out-Null
$d.Content.Find.Execute("57+32=89", $true, $false, $false, $false, $false, $true, 1, $false, "46+34=80", 2) | Out-Null
$d.Content.Find.Execute("41+32=73", $true, $false, $false, $false, $false, $true, 1, $false, "44+44=88", 2) | Out-Null
$d.Content.Find.Execute("53+17=70", $true, $false, $false, $false, $false, $true, 1, $false, "35+46=81", 2) | Out-Null
$d.Content.Find.Execute("86-40=46", $true, $false, $false, $false, $false, $true, 1, $false, "46-18=28", 2) | Out-Null
$d.Content.Find.Execute("30-9=21", $true, $false, $false, $false, $false, $true, 1, $false, "0+68=68", 2) | Out-Null
$d.Content.Find.Execute("54-15=39", $true, $false, $false, $false, $false, $true, 1, $false, "10+31=41", 2) | Out-Null
$d.Content.Find.Execute("78-46=32", $true, $false, $false, $false, $false, $true, 1, $false, "11-0=11", 2) | Out-Null
$d.Content.Find.Execute("23+64=87", $true, $false, $false, $false, $false, $true, 1, $false, "20+38=58", 2) | Out-Null
$d.Content.Find.Execute("31-13=18", $true, $false, $false, $false, $false, $true, 1, $false, "19+73=92", 2) | Out-Null
$d.Content.Find.Execute("83-30=53", $true, $false, $false, $false, $false, $true, 1, $false, "19+35=54", 2) | Out-Null
$d.Content.Find.Execute("95-53=42", $true, $false, $false, $false, $false, $true, 1, $false, "72-12=60", 2) | Out-Null
$d.Content.Find.Execute("3+17=20", $true, $false, $false, $false, $false, $true, 1, $false, "66+15=81", 2) | Out-Null
$d.Content.Find.Execute("61-9=52", $true, $false, $false, $false, $false, $true, 1, $false, "22+62=84", 2) | Out-Null
$d.Content.Find.Execute("50+16=66", $true, $false, $false, $false, $false, $true, 1, $false, "49+27=76", 2) | Out-Null
$d.Content.Find.Execute("20+69=89", $true, $false, $false, $false, $false, $true, 1, $false, "60-29=31", 2) | Out-Null
$d.Content.Find.Execute("99-25=74", $true, $false, $false, $false, $false, $true, 1, $false, "4+66=70", 2) | Out-Null
$d.Content.Find.Execute("15+30=45", $true, $false, $false, $false, $false, $true, 1, $false, "15+71=86", 2) | Out-Null
$d.Content.Find.Execute("65-7=58", $true, $false, $false, $false, $false, $true, 1, $false, "52+33=85", 2) | Out-Null
$d.Content.Find.Execute("48+22=70", $true, $false, $false, $false, $false, $true, 1, $false, "26+53=79", 2) | Out-Null
$d.Content.Find.Execute("93-29=64", $true, $false, $false, $false, $false, $true, 1, $false, "5+67=72", 2) | Out-Null
$d.Content.Find.Execute("60-53=7", $true, $false, $false, $false, $false, $true, 1, $false, "41+49=90", 2) | Out-Null
$d.Content.Find.Execute("78-62=16", $true, $false, $false, $false, $false, $true, 1, $false, "44+51=95", 2) | Out-Null
$d.Content.Find.Execute("86-82=4", $true, $false, $false, $false, $false, $true, 1, $false, "97-37=60", 2) | Out-Null
$d.Content.Find.Execute("83-79=4", $true, $false, $false, $false, $false, $true, 1, $false, "58-31=27", 2) | Out-Null
$d.Content.Find.Execute("12+29=41", $true, $false, $false, $false, $false, $true, 1, $false, "55-32=23", 2) | Out-Null
$d.Content.Find.Execute("41+55=96", $true, $false, $false, $false, $false, $true, 1, $false, "90-29=61", 2) | Out-Null
$d.Content.Find.Execute("85-29=56", $true, $false, $false, $false, $false, $true, 1, $false, "22+54=76", 2) | Out-Null
